# Applies the cryptos-list price/volume refresh for Sat Aug 10 05:42:51 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every written cell to be stored as literal text (matches the original
# inlineStr cells) instead of letting Excel auto-coerce numeric-looking strings
# (e.g. "509.20", "0.0234") into floating point numbers, which would both change
# the stored type and silently drop significant trailing/representational digits.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.293.77"
Set-TextValue "E2" "  -1.06%  "
Set-TextValue "D3" "2.593.47"
Set-TextValue "E3" "  -2.07%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "509.20"
Set-TextValue "E5" "  -0.49%  "
Set-TextValue "D6" "153.64"
Set-TextValue "E6" "  -2.42%  "
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "E8" "  -3.03%  "
Set-TextValue "D9" "2.600.23"
Set-TextValue "E9" "  -2.60%  "
Set-TextValue "D10" "6.67"
Set-TextValue "E10" "  +4.66%  "
Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  -1.34%  "
Set-TextValue "E12" "  -0.92%  "
Set-TextValue "E13" "  +1.64%  "
Set-TextValue "D14" "3.048.56"
Set-TextValue "E14" "  -2.42%  "
Set-TextValue "D15" "60.270.83"
Set-TextValue "E15" "  -1.10%  "
Set-TextValue "E16" "  -1.49%  "
Set-TextValue "E17" "  +0.09%  "
Set-TextValue "D18" "2.602.10"
Set-TextValue "E18" "  -2.41%  "
Set-TextValue "D19" "4.73"
Set-TextValue "E19" "  -1.40%  "
Set-TextValue "D20" "352.80"
Set-TextValue "E20" "  +1.17%  "
Set-TextValue "D21" "10.51"
Set-TextValue "E21" "  +0.03%  "
Set-TextValue "D22" "6.13"
Set-TextValue "E22" "  -0.85%  "
Set-TextValue "E23" "  +0.00%  "
Set-TextValue "E24" "  +0.16%  "
Set-TextValue "D25" "0.419"
Set-TextValue "E25" "  -0.59%  "
Set-TextValue "E26" "  -0.41%  "
Set-TextValue "E27" "  +0.21%  "
Set-TextValue "D28" "0.0₃0836"
Set-TextValue "E28" "  -2.90%  "
Set-TextValue "E29" "  -2.60%  "
Set-TextValue "E30" "  -0.08%  "
Set-TextValue "D31" "19.32"
Set-TextValue "E31" "  -0.94%  "
Set-TextValue "D32" "150.97"
Set-TextValue "E32" "  -4.32%  "
Set-TextValue "E33" "  -1.04%  "
Set-TextValue "D34" "5.73"
Set-TextValue "E34" "  +0.25%  "
Set-TextValue "D35" "3.98"
Set-TextValue "E35" "  -1.61%  "
Set-TextValue "E36" "  -3.14%  "
Set-TextValue "D37" "0.873"
Set-TextValue "E37" "  +4.21%  "
Set-TextValue "E38" "  -3.20%  "
Set-TextValue "D39" "36.23"
Set-TextValue "E39" "  +2.27%  "
Set-TextValue "E40" "  -2.55%  "
Set-TextValue "E41" "  -0.55%  "
Set-TextValue "D42" "294.42"
Set-TextValue "E42" "  -5.01%  "
Set-TextValue "E43" "  -0.49%  "
Set-TextValue "D44" "0.617"
Set-TextValue "E44" "  -4.03%  "
Set-TextValue "E45" "  +0.07%  "
Set-TextValue "D46" "0.0553"
Set-TextValue "E46" "  -4.42%  "
Set-TextValue "D47" "19.59"
Set-TextValue "E47" "  -1.45%  "
Set-TextValue "B48" "VeChain"
Set-TextValue "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D48" "0.0234"
Set-TextValue "E48" "  -1.22%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "4.79"
Set-TextValue "E49" "  -1.31%  "
Set-TextValue "D50" "10.32"
Set-TextValue "E50" "  -0.11%  "
Set-TextValue "D51" "1.987.23"
Set-TextValue "E51" "  -2.76%  "
